$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G (shifts the existing "FK - AccountID" column from G to H,
# in the Receipts table header/data rows 7-8) and give it a title + sample value.
$ws.Columns("G:G").Insert()
$ws.Range("G7").Value = "Contributions"
$ws.Range("G8").Value = "Dan,£5/Marie,£5"

# Best-fit the new column's width, same as Excel auto-sizing it to the text it holds.
$ws.Columns("G:G").AutoFit()

# Match the author's final view state: scrolled right a bit with G15 selected.
[void]$ws.Range("G15").Select()

"done"
